# Feb 2 - Cart Language Test Data update
# Remove obsolete test-data rows from Sheet1 (rows 7, 8, 10, 18 in original numbering).
# Delete from the bottom up so earlier row numbers stay valid as we go.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$rowsToDelete = @(18, 10, 8, 7)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
